$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 658
$wsExpo.Range("F6").Value = 1586
$wsExpo.Range("F8").Value = 3152

# Sheet "全部类型" (all types) - same events duplicated, but at different rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 658
$wsAll.Range("F7").Value = 1586
$wsAll.Range("F9").Value = 3152
